$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.902.58"
$ws.Range("E2").Value = "  +0.15%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.546.69"
$ws.Range("E3").Value = "  -1.03%  "
$ws.Range("E4").Value = "  +0.33%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "206.49"
$ws.Range("E5").Value = "  +0.30%  "
$ws.Range("E6").Value = "  -0.24%  "
$ws.Range("E7").Value = "  +0.29%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "21.33"
$ws.Range("E9").Value = "  -1.67%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0583"
$ws.Range("E10").Value = "  -0.12%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0857"
$ws.Range("E11").Value = "  -1.03%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.766.79"
$ws.Range("E12").Value = "  -1.05%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.547.99"
$ws.Range("E13").Value = "  -1.22%  "
$ws.Range("E14").Value = "  -0.65%  "
$ws.Range("E15").Value = "  -0.41%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "26.899.97"
$ws.Range("E16").Value = "  +0.06%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.43"
$ws.Range("E17").Value = "  +0.29%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "215.66"
$ws.Range("E18").Value = "  +0.46%  "
$ws.Range("E19").Value = "  +0.55%  "
$ws.Range("E20").Value = "  -2.04%  "
$ws.Range("E21").Value = "  +0.35%  "
$ws.Range("E22").Value = "  -2.55%  "
$ws.Range("E23").Value = "  +0.23%  "
$ws.Range("E24").Value = "  -3.08%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.98"
$ws.Range("E25").Value = "  -1.28%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.61"
$ws.Range("E26").Value = "  -1.99%  "
$ws.Range("E27").Value = "  -0.62%  "
$ws.Range("E28").Value = "  +0.29%  "
$ws.Range("E29").Value = "  +0.27%  "
$ws.Range("E30").Value = "  -1.02%  "
$ws.Range("E31").Value = "  -1.07%  "
$ws.Range("E32").Value = "  +1.74%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.358.68"
$ws.Range("E33").Value = "  -3.14%  "
$ws.Range("E34").Value = "  +0.77%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.961"
$ws.Range("E36").Value = "  +4.77%  "
$ws.Range("E37").Value = "  +0.20%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0164"
$ws.Range("E38").Value = "  -0.12%  "
$ws.Range("E39").Value = "  -1.05%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.805"
$ws.Range("E40").Value = "  -0.79%  "
$ws.Range("E41").Value = "  +0.29%  "
$ws.Range("E42").Value = "  +3.81%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.989"
$ws.Range("E43").Value = "  -0.75%  "
$ws.Range("E44").Value = "  +2.09%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "63.40"
$ws.Range("E45").Value = "  +0.37%  "
$ws.Range("E46").Value = "  -1.82%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.680.88"
$ws.Range("E47").Value = "  -1.18%  "
$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "85.74"
$ws.Range("E48").Value = "  -0.62%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0511"
$ws.Range("E49").Value = "  +1.06%  "
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0₇0971"
$ws.Range("E50").Value = "  -0.98%  "
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0947"
$ws.Range("E51").Value = "  +0.05%  "
